# "Generate Report for Handback" — records that the de-de handback just
# completed (new target/handback files + handback timestamps) and that the
# zh-cn handback (already recorded) is now reflected as "in sync" in the
# overview/status columns.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status text refresh (was "Ready for handoff") ------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- URLs for the two source docs (same as the existing "A" hyperlinks) ---
$url31952 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1047de76610e913e7e45ea06f340eead55fb247e/e2e/31952d01-eb04-4548-951c-564b77d6041b.md"
$url42a3f = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1047de76610e913e7e45ea06f340eead55fb247e/e2e/42a3fa32-9b3e-45e6-ad91-1ab53e63cf42.md"

$name31952 = "31952d01-eb04-4548-951c-564b77d6041b.md"
$name42a3f = "42a3fa32-9b3e-45e6-ad91-1ab53e63cf42.md"

# --- zh-cn: Latest Target File (I) + Latest Handback File (J) -------------
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url31952, "", "", $name31952)
$wsZhCn.Range("J2").Value = "31952d01-eb04-4548-951c-564b77d6041b.1cbbf510475b4b817a2d09c21943873666a77504.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-24 16:29:39"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url42a3f, "", "", $name42a3f)
$wsZhCn.Range("J3").Value = "42a3fa32-9b3e-45e6-ad91-1ab53e63cf42.fb5a7979b422cf48086abdafb311712ca845fceb.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-24 16:29:39"

# --- de-de: Latest Target File (I) + Latest Handback File (J/K) -----------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url31952, "", "", $name31952)
$wsDeDe.Range("J2").Value = "31952d01-eb04-4548-951c-564b77d6041b.1cbbf510475b4b817a2d09c21943873666a77504.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-24 16:29:46"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url42a3f, "", "", $name42a3f)
$wsDeDe.Range("J3").Value = "42a3fa32-9b3e-45e6-ad91-1ab53e63cf42.fb5a7979b422cf48086abdafb311712ca845fceb.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-24 16:29:46"

Write-Host "Handback report generated."
